$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 75727.71000000001
$ws.Range("I21").Value = 88019
$ws.Range("J21").Value = 44999.5
$ws.Range("K21").Value = 88019
$ws.Range("L21").Value = 44999.5
$ws.Range("M21").Value = -87551
$ws.Range("N21").Value = -45935.5
$ws.Range("H23").Value = 75727.71000000001
$ws.Range("I23").Value = 88019
$ws.Range("J23").Value = 44999.5
$ws.Range("K23").Value = 88019
$ws.Range("L23").Value = 44999.5
$ws.Range("M23").Value = -87785
$ws.Range("N23").Value = -45467.5
$ws.Range("H64").Value = 2876.8462
$ws.Range("I64").Value = 2929.9
$ws.Range("J64").Value = 2700
$ws.Range("K64").Value = 2929.9
$ws.Range("L64").Value = 2700
$ws.Range("M64").Value = -2681.9
$ws.Range("N64").Value = -3196
$ws.Range("H67").Value = 2876.8462
$ws.Range("I67").Value = 2929.9
$ws.Range("J67").Value = 2700
$ws.Range("K67").Value = 2929.9
$ws.Range("L67").Value = 2700
$ws.Range("M67").Value = -2071.9
$ws.Range("N67").Value = -4416
$ws.Range("H74").Value = 11114501
$ws.Range("I74").Value = 11114501
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 11114501
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -11113565
$ws.Range("N74").ClearContents()
$ws.Range("H76").Value = 3633.3333
$ws.Range("I76").Value = 3300
$ws.Range("J76").Value = 3800
$ws.Range("K76").Value = 3300
$ws.Range("L76").Value = 3800
$ws.Range("N76").Value = -4430
$ws.Range("M76").Value = -2985
$ws.Range("H77").Value = 11114501
$ws.Range("I77").Value = 11114501
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 55572505
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -55567825
$ws.Range("N77").ClearContents()
$ws.Range("H79").Value = 3633.3333
$ws.Range("I79").Value = 3300
$ws.Range("J79").Value = 3800
$ws.Range("K79").Value = 3300
$ws.Range("L79").Value = 3800
$ws.Range("N79").Value = -5984
$ws.Range("M79").Value = -2208
$ws.Range("H112").Value = 1593.84
$ws.Range("J112").Value = 1609.1237
$ws.Range("L112").Value = 4827.3711
$ws.Range("N112").Value = -7043.3711
$ws.Range("H128").Value = 42580
$ws.Range("J128").Value = 42580
$ws.Range("L128").Value = 42580
$ws.Range("N128").Value = -52540
$ws.Range("H130").Value = 42413.332
$ws.Range("J130").Value = 42413.332
$ws.Range("L130").Value = 42413.332
$ws.Range("N130").Value = -52453.332
$ws.Range("H133").Value = 52867
$ws.Range("J133").Value = 52867
$ws.Range("L133").Value = 52867
$ws.Range("N133").Value = -62987
$ws.Range("H138").Value = 2369.65
$ws.Range("I138").Value = 1200.8
$ws.Range("J138").Value = 2999.0308
$ws.Range("K138").Value = 3602.4
$ws.Range("L138").Value = 8997.0924
$ws.Range("M138").Value = 1537.6
$ws.Range("N138").Value = -19277.0924

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 31254708
$ws.Range("I62").Value = 55559710
$ws.Range("J62").Value = 5420
$ws.Range("K62").Value = 55559710
$ws.Range("L62").Value = 5420
$ws.Range("M62").Value = -55559086
$ws.Range("N62").Value = -6668
$ws.Range("H65").Value = 31254708
$ws.Range("I65").Value = 55559710
$ws.Range("J65").Value = 5420
$ws.Range("K65").Value = 277798550
$ws.Range("L65").Value = 27100
$ws.Range("M65").Value = -277795430
$ws.Range("N65").Value = -33340
$ws.Range("H88").Value = 24945
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 29926.666
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 29926.666
$ws.Range("M88").Value = -9594
$ws.Range("N88").Value = -30738.666
$ws.Range("H91").Value = 24945
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 29926.666
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 29926.666
$ws.Range("M91").Value = -8596
$ws.Range("N91").Value = -32734.666
$ws.Range("H97").Value = 34600
$ws.Range("J97").Value = 34600
$ws.Range("L97").Value = 34600
$ws.Range("N97").Value = -36582
$ws.Range("H99").Value = 4082.0908
$ws.Range("I99").Value = 950.5
$ws.Range("J99").Value = 5871.5713
$ws.Range("K99").Value = 950.5
$ws.Range("L99").Value = 5871.5713
$ws.Range("M99").Value = 547.5
$ws.Range("N99").Value = -8867.5713
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H112").Value = 30178.572
$ws.Range("J112").Value = 30178.572
$ws.Range("L112").Value = 30178.572
$ws.Range("N112").Value = -33132.572
$ws.Range("H118").Value = 28390
$ws.Range("J118").Value = 28390
$ws.Range("L118").Value = 28390
$ws.Range("N118").Value = -31704
$ws.Range("H122").Value = 2448.04
$ws.Range("I122").Value = 1868.2941
$ws.Range("J122").Value = 3680
$ws.Range("K122").Value = 5604.8823
$ws.Range("L122").Value = 11040
$ws.Range("M122").Value = -3154.8823
$ws.Range("N122").Value = -15940
$ws.Range("H126").Value = 4082.0908
$ws.Range("I126").Value = 950.5
$ws.Range("J126").Value = 5871.5713
$ws.Range("K126").Value = 2851.5
$ws.Range("L126").Value = 17614.7139
$ws.Range("M126").Value = -381.5
$ws.Range("N126").Value = -22554.7139

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 14962.4
$ws.Range("I34").Value = 36916
$ws.Range("J34").Value = 9474
$ws.Range("K34").Value = 110748
$ws.Range("L34").Value = 28422
$ws.Range("M34").Value = -110664
$ws.Range("N34").Value = -28590
$ws.Range("H39").Value = 7267.8335
$ws.Range("J39").Value = 7530.5884
$ws.Range("L39").Value = 22591.7652
$ws.Range("N39").Value = -23179.7652
$ws.Range("H55").Value = 203202
$ws.Range("I55").Value = 501000
$ws.Range("J55").Value = 4670
$ws.Range("K55").Value = 1503000
$ws.Range("L55").Value = 14010
$ws.Range("M55").Value = -1502823
$ws.Range("N55").Value = -14364
$ws.Range("H121").Value = 2256.4707
$ws.Range("J121").Value = 2420.2341
$ws.Range("L121").Value = 7260.702300000001
$ws.Range("N121").Value = -9880.702300000001
$ws.Range("H122").Value = 2462.5874
$ws.Range("J122").Value = 3506.2195
$ws.Range("L122").Value = 31555.9755
$ws.Range("N122").Value = -36455.9755
$ws.Range("H123").Value = 3749.75
$ws.Range("I123").Value = 3500
$ws.Range("J123").Value = 3999.5
$ws.Range("K123").Value = 10500
$ws.Range("L123").Value = 11998.5
$ws.Range("M123").Value = -8050
$ws.Range("N123").Value = -16898.5
$ws.Range("H124").Value = 6500
$ws.Range("J124").Value = 6500
$ws.Range("L124").Value = 19500
$ws.Range("N124").Value = -29320
$ws.Range("H125").Value = 8200
$ws.Range("J125").Value = 10000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -39840

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2248.743
$ws.Range("I132").Value = 1104.3
$ws.Range("K132").Value = 3312.9
$ws.Range("M132").Value = -782.8999999999996
$ws.Range("H133").Value = 58016
$ws.Range("J133").Value = 56433.332
$ws.Range("L133").Value = 56433.332
$ws.Range("N133").Value = -66553.33199999999
$ws.Range("H135").Value = 62002.8
$ws.Range("J135").Value = 56467.285
$ws.Range("L135").Value = 56467.285
$ws.Range("N135").Value = -66607.285
$ws.Range("H140").Value = 42750
$ws.Range("J140").Value = 42750
$ws.Range("L140").Value = 42750
$ws.Range("N140").Value = -53110

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 38200
$ws.Range("J111").Value = 38200
$ws.Range("L111").Value = 38200
$ws.Range("N111").Value = -46380

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 16252.4
$ws.Range("J54").Value = 16252.4
$ws.Range("L54").Value = 16252.4
$ws.Range("N54").Value = -17292.4
